# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Papaya" (Vega Modelo de Temuco)
# just above the existing row 123, pushing the subsequent rows down by one
# and extending the sheet from A1:T142 to A1:T143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 123 (shifts rows 123:142 down to 124:143,
# carrying their formatting with them).
$ws.Range("A123").EntireRow.Insert()

# Populate the newly inserted row 123 with the new observation.
$ws.Range("A123").Value = 10
$ws.Range("B123").Value = 'Vega Modelo de Temuco'
$ws.Range("C123").Value = 'La Araucanía'
$ws.Range("D123").Value = 45258
$ws.Range("E123").Value = 9
$ws.Range("F123").Value = 'Fruta'
$ws.Range("G123").Value = 100108
$ws.Range("H123").Value = 'Tropicales y subtropicales'
$ws.Range("I123").Value = 100108004
$ws.Range("J123").Value = 'Papaya'
$ws.Range("K123").Value = 'Cultivar IV Región'
$ws.Range("L123").Value = 'Primera'
$ws.Range("M123").Value = 30
$ws.Range("N123").Value = 2500
$ws.Range("O123").Value = 2500
$ws.Range("P123").Value = 2500
$ws.Range("Q123").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R123").Value = 'Provincia del Elquí'
$ws.Range("S123").Value = 2500
$ws.Range("T123").Value = 1
